$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-26 09:37:31"
$wsZh.Range("G2").Value = "2016-01-26 09:38:19"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-26 09:37:44"
$wsDe.Range("G2").Value = "2016-01-26 09:38:41"
